{"js": "// The \"Requisitos\" section of this course-catalog page lists prerequisite\n// disciplines as a sequence of lines (one requisite per line, separated by\n// manual line breaks) inside a single List-Bullet paragraph. The edit\n// re-orders those lines (same 24 requisites, new order) without touching\n// anything else in the document.\n\n// New position -> old position (0-based) inside the requisites list.\nconst NEW_ORDER = [\n  22, 19, 3, 2, 5, 10, 12, 1, 8, 14, 16, 23, 15, 17, 21, 7, 18, 4, 6, 13, 9, 20, 11, 0\n];\n\nfunction escapeXml(str) {\n  return str\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Find the \"Requisitos\" heading, then the very next paragraph (the bullet\n// list that actually holds the requisite lines).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text,style\"));\nawait context.sync();\n\nconst heading = paragraphs.items.find(\n  (p) => p.style === \"Heading 2\" && p.text.trim() === \"Requisitos\"\n);\nif (!heading) {\n  throw new Error('Could not find the \"Requisitos\" heading paragraph.');\n}\n\nconst list = heading.getNext();\nlist.load(\"text,style\");\nawait context.sync();\n\n// Each requisite is on its own line, separated by the manual line breaks\n// (represented as vertical-tab \"\\v\" in the Word text model). The source\n// paragraph ends with a trailing break, so splitting on \"\\v\" yields one\n// empty trailing element that we drop.\nconst lines = list.text.split(\"\\v\");\nif (lines.length && lines[lines.length - 1] === \"\") {\n  lines.pop();\n}\n\nif (lines.length !== NEW_ORDER.length) {\n  throw new Error(\n    `Expected ${NEW_ORDER.length} requisite lines, found ${lines.length}.`\n  );\n}\n\nconst reordered = NEW_ORDER.map((oldIndex) => lines[oldIndex]);\n\n// Rebuild the paragraph: one run per requisite line, each run holding its\n// own text plus the trailing line break - matching the original structure.\nconst runsXml = reordered\n  .map((line) => `<w:r><w:t>${escapeXml(line)}</w:t><w:br/></w:r>`)\n  .join(\"\");\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>\n            ${runsXml}\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst range = list.getRange();\nrange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The \"Requisitos\" section of this course-catalog page lists prerequisite\n# disciplines as a sequence of lines (one requisite per line, separated by\n# manual line breaks) inside a single List-Bullet paragraph. The edit\n# re-orders those lines (same 24 requisites, new order) without touching\n# anything else in the document.\n\n$doc = $word.ActiveDocument\n$paras = $doc.Paragraphs\n\n# Locate the \"Requisitos\" heading paragraph, then the very next paragraph\n# (the bullet list that actually holds the requisite lines).\n$heading = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $candidate = $paras.Item($i)\n    if ($candidate.Style.NameLocal -eq \"Heading 2\" -and $candidate.Range.Text.Trim() -eq \"Requisitos\") {\n        $heading = $candidate\n        break\n    }\n}\nif ($heading -eq $null) {\n    throw \"Could not find the 'Requisitos' heading paragraph.\"\n}\n\n$list = $heading.Next()\n\n# Each requisite is on its own line, separated by manual line breaks\n# (represented as vertical-tab \"`v\" in the Word text model). Range.Text\n# also carries a trailing paragraph-mark character (CR) at the very end,\n# and the source paragraph ends with a trailing break, so after trimming\n# the CR and splitting on \"`v\" we get one empty trailing element to drop.\n$fullText = $list.Range.Text.TrimEnd([char]13)\n$lines = $fullText -split \"`v\"\nif ($lines.Count -gt 0 -and $lines[$lines.Count - 1] -eq \"\") {\n    $lines = $lines[0..($lines.Count - 2)]\n}\n\n# New position -> old position (0-based) inside the requisites list.\n$newOrder = @(22, 19, 3, 2, 5, 10, 12, 1, 8, 14, 16, 23, 15, 17, 21, 7, 18, 4, 6, 13, 9, 20, 11, 0)\n\nif ($lines.Count -ne $newOrder.Count) {\n    throw \"Expected $($newOrder.Count) requisite lines, found $($lines.Count).\"\n}\n\n$reordered = @()\nforeach ($oldIndex in $newOrder) {\n    $reordered += $lines[$oldIndex]\n}\n\nfunction Escape-Xml([string]$text) {\n    return $text.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\").Replace('\"', \"&quot;\")\n}\n\n# Rebuild the paragraph: one run per requisite line, each run holding its\n# own text plus the trailing line break - matching the original structure.\n$runsXml = \"\"\nforeach ($line in $reordered) {\n    $runsXml += \"<w:r><w:t>\" + (Escape-Xml $line) + \"</w:t><w:br/></w:r>\"\n}\n\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>\n            $runsXml\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$start = $list.Range.Start\n$end = $list.Range.End\n$fullRange = $doc.Range($start, $end)\n$fullRange.InsertXML($ooxml)\n"}
